$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing users (user1 -> u1, user2 -> u2), shorten their passwords
$ws.Range("A2").Value = "u1"
$ws.Range("B2").Value = 11
$ws.Range("A3").Value = "u2"
$ws.Range("B3").Value = 22

# Add four new user/password rows
$ws.Range("A4").Value = "u3"
$ws.Range("B4").Value = 33
$ws.Range("A5").Value = "u4"
$ws.Range("B5").Value = 44
$ws.Range("A6").Value = "u5"
$ws.Range("B6").Value = 55
$ws.Range("A7").Value = "u6"
$ws.Range("B7").Value = 66

# Match the saved selection from the edited workbook
$ws.Range("E9").Select()
